$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(0.237, 0.192, 0.275, 0.258, 0.111, 0.128, 0.128, 0.148, 0.158, 0.218)
    3 = @(0.293, 0.264, 0.341, 0.305, 0.126, 0.2,   0.207, 0.207, 0.215, 0.349)
    4 = @(0.438, 0.416, 0.466, 0.434, 0.15,  0.404, 0.402, 0.408, 0.454, 0.659)
    5 = @(0.572, 0.539, 0.591, 0.584, 0.211, 0.626, 0.582, 0.594, 0.649, 0.838)
    6 = @(0.671, 0.66,  0.6830000000000001, 0.639, 0.235, 0.755, 0.724, 0.698, 0.792, 0.9330000000000001)
    7 = @(0.74,  0.752, 0.775, 0.756, 0.275, 0.87,  0.822, 0.8100000000000001, 0.879, 0.964)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($i = 0; $i -lt $values.Count; $i++) {
        # Column B is index 2, so offset by 2
        $col = $i + 2
        $ws.Cells.Item($row, $col).Value = $values[$i]
    }
}

# Remove row 8 entirely (the "50" row), shrinking the used range to A1:L7
$ws.Rows.Item(8).Delete() | Out-Null
